$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.915.01"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.359.14"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.87"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.42"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.95"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.936.29"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.96"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.356.65"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.004.97"
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.35"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.20"
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.82"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.86"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.558"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.500.99"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("E26").Value = "  -6.95%  "
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("E32").Value = "  -4.66%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "170.24"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.77"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.99"
$ws.Range("E39").Value = "  -8.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.394.82"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0743"
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.31"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.487.52"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.43"
$ws.Range("E50").Value = "  -1.21%  "
